$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# Add the new row of parameters: fix_prices_to_2020
$ws.Range("A15").Value = "fix_prices_to_2020"
$ws.Range("B15").Value = $true
$ws.Range("C15").Value = "for verification runs. Fix fuel prices, CO2 and electricity demand to 2020 data"

# Update the selected cell on this sheet
$ws.Range("C9").Select()
